$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (attendee count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1900
$wsExpo.Range("F4").Value = 825
$wsExpo.Range("F5").Value = 816
$wsExpo.Range("F6").Value = 267

# Sheet "全部类型" - update 想去人数 (attendee count) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1900
$wsAll.Range("F5").Value = 825
$wsAll.Range("F6").Value = 816
$wsAll.Range("F7").Value = 267
